# Update sheet1 with new TPM values and drop the self-pair rows (old rows 8-10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old rows 8, 9 and 10 (self-pairs for MuSCs sending cluster) - entire rows
$ws.Rows.Item(8).Resize(3).Delete()

# Row 2 (ECs -> FAPs)
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 0.02027766666666667
$ws.Range("N2").Value = 0.060833
$ws.Range("O2").Value = 0.7477567175553752
$ws.Range("P2").Value = 0.7477567175553753
$ws.Range("Q2").Value = 1.682796532757667
$ws.Range("R2").Value = 15.145168794819
$ws.Range("S2").Value = 0.3357056860803754
$ws.Range("T2").Value = 0.3357056860803754

# Row 3 (ECs -> MuSCs)
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 82.98768099999999
$ws.Range("H3").Value = 248.963043
$ws.Range("I3").Value = 0.4489504115427952
$ws.Range("J3").Value = 0.4489504115427952
$ws.Range("M3").Value = 0.006840333333333334
$ws.Range("N3").Value = 0.020521
$ws.Range("O3").Value = 0.2522432824446247
$ws.Range("P3").Value = 0.2522432824446247
$ws.Range("Q3").Value = 0.5676634006003334
$ws.Range("R3").Value = 5.108970605403
$ws.Range("S3").Value = 0.1132447254624198
$ws.Range("T3").Value = 0.1132447254624198

# Row 4 (FAPs -> FAPs)
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 63.14058933333333
$ws.Range("H4").Value = 189.421768
$ws.Range("I4").Value = 0.3415807409566563
$ws.Range("J4").Value = 0.3415807409566563
$ws.Range("M4").Value = 0.02027766666666667
$ws.Range("N4").Value = 0.060833
$ws.Range("O4").Value = 0.7477567175553752
$ws.Range("P4").Value = 0.7477567175553753
$ws.Range("Q4").Value = 1.280343823638222
$ws.Range("R4").Value = 11.523094412744
$ws.Range("S4").Value = 0.2554192936378822
$ws.Range("T4").Value = 0.2554192936378823

# Row 5 (FAPs -> MuSCs)
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("I5").Value = 0.3415807409566563
$ws.Range("J5").Value = 0.3415807409566563
$ws.Range("M5").Value = 0.006840333333333334
$ws.Range("N5").Value = 0.020521
$ws.Range("O5").Value = 0.2522432824446247
$ws.Range("P5").Value = 0.2522432824446247
$ws.Range("Q5").Value = 0.4319026779031111
$ws.Range("R5").Value = 3.887124101128
$ws.Range("S5").Value = 0.08616144731877405
$ws.Range("T5").Value = 0.08616144731877405

# Row 6 (MuSCs -> FAPs)
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 38.719942
$ws.Range("H6").Value = 116.159826
$ws.Range("I6").Value = 0.2094688475005485
$ws.Range("J6").Value = 0.2094688475005485
$ws.Range("N6").Value = 0.060833
$ws.Range("O6").Value = 0.7477567175553752
$ws.Range("P6").Value = 0.7477567175553753
$ws.Range("Q6").Value = 0.7851500772286667
$ws.Range("R6").Value = 7.066350695058
$ws.Range("S6").Value = 0.1566317378371176
$ws.Range("T6").Value = 0.1566317378371176

# Row 7 (MuSCs -> MuSCs)
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 38.719942
$ws.Range("H7").Value = 116.159826
$ws.Range("I7").Value = 0.2094688475005485
$ws.Range("J7").Value = 0.2094688475005485
$ws.Range("O7").Value = 0.2522432824446247
$ws.Range("P7").Value = 0.2522432824446247
$ws.Range("Q7").Value = 0.2648573099273334
$ws.Range("R7").Value = 2.383715789346
$ws.Range("S7").Value = 0.05283710966343088
$ws.Range("T7").Value = 0.05283710966343088
